$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current values (rows 2..18, columns A..E = 1..5) using .Text
# (the .Value getter on empty cells misbehaves in this runtime, .Text is reliable
# and all data in this sheet is plain text / shared strings).
$data = @()
for ($r = 2; $r -le 18; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 5; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Text
    }
    $data += ,$rowVals
}

# $data[0] holds old row 2 ... $data[16] holds old row 18.
# New order: new row2 = old row18; new rows 3..18 = old rows 2..17 (shifted down by one).
$oldRow18 = $data[16]
$shifted = $data[0..15]

$newOrder = @()
$newOrder += ,$oldRow18
$newOrder += $shifted

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $r = $i + 2
    $rowVals = $newOrder[$i]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Match the author's final selection/scroll state: whole row 19 selected.
$win = $excel.ActiveWindow
$win.Zoom = 223
$win.TopLeftCell = $ws.Range("A9")
$ws.Rows("19:19").EntireRow.Select() | Out-Null
